$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 940
$ws.Range("F4").Value = 52
$ws.Range("F5").Value = 227
$ws.Range("F6").Value = 47
$ws.Range("F7").Value = 1168
$ws.Range("F8").Value = 933
$ws.Range("F9").Value = 29
$ws.Range("F10").Value = 724
$ws.Range("F12").Value = 1483
$ws.Range("F13").Value = 64
$ws.Range("F14").Value = 142
$ws.Range("F15").Value = 1649
$ws.Range("F17").Value = 627
$ws.Range("F18").Value = 18
$ws.Range("F19").Value = 10
$ws.Range("F21").Value = 1090
$ws.Range("F22").Value = 1516
$ws.Range("F23").Value = 759
$ws.Range("F24").Value = 633
$ws.Range("F25").Value = 506
$ws.Range("F26").Value = 478
$ws.Range("F28").Value = 14
$ws.Range("F30").Value = 1155
$ws.Range("F31").Value = 318
$ws.Range("F32").Value = 2439
$ws.Range("F34").Value = 1392
$ws.Range("F35").Value = 465
$ws.Range("F38").Value = 4001
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 184
$ws.Range("F14").Value = 4138
$ws.Range("F22").Value = 260
$ws.Range("F30").Value = 1718
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 98
$ws.Range("F5").Value = 1673
$ws.Range("F7").Value = 1021
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1673
$ws.Range("F6").Value = 1021
$ws.Range("F8").Value = 940
$ws.Range("F9").Value = 52
$ws.Range("F10").Value = 47
$ws.Range("F11").Value = 1168
$ws.Range("F12").Value = 933
$ws.Range("F14").Value = 29
$ws.Range("F16").Value = 724
$ws.Range("F17").Value = 184
$ws.Range("F18").Value = 184
$ws.Range("F21").Value = 1483
$ws.Range("F22").Value = 64
$ws.Range("F23").Value = 142
$ws.Range("F24").Value = 1649
$ws.Range("F26").Value = 627
$ws.Range("F29").Value = 1090
$ws.Range("F30").Value = 1516
$ws.Range("F31").Value = 759
$ws.Range("F32").Value = 633
$ws.Range("F33").Value = 506
$ws.Range("F34").Value = 478
$ws.Range("F37").Value = 260
$ws.Range("F40").Value = 1155
$ws.Range("F41").Value = 318
$ws.Range("F42").Value = 2439
$ws.Range("F45").Value = 1718
$ws.Range("F46").Value = 1718
$ws.Range("F47").Value = 1392
$ws.Range("F48").Value = 465
$ws.Range("F50").Value = 4001
